$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @("2.900324425070266e-11", "0.0001537489499301437", "1935279062.313128", "198602002.3250627", "0", "2133881064.638345")
    3 = @("0.3464964993005633", "0.3375848360084654", "0.7127328510149897", "6.48142807727062", "0", "7.878242263594639")
    4 = @("0.3464964993005633", "1.65323645889881", "3.082599426703578", "6.48142807727062", "0", "11.56376046217357")
    5 = @("1.505614041169197", "1.65323645889881", "3.082599426703578", "246.9852506941017", "0", "253.2267006208733")
    6 = @("1.505614041169197", "9.226618575922256", "0.7127328510149897", "6.48142807727062", "0", "17.92639354537706")
    7 = @("0.06328177979961902", "0.3375848360084654", "0.7127328510149897", "6.48142807727062", "0", "7.595027544093695")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = [double]$vals[0]
    $ws.Cells.Item($row, 3).Value = [double]$vals[1]
    $ws.Cells.Item($row, 4).Value = [double]$vals[2]
    $ws.Cells.Item($row, 5).Value = [double]$vals[3]
    $ws.Cells.Item($row, 6).Value = [double]$vals[4]
    $ws.Cells.Item($row, 7).Value = [double]$vals[5]
}
